$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: requirement "Prototipo de Interfaz web" had its status changed
# from "Viernes 14" to "caca"
$ws.Range("I30").Value = "caca"

# New requisito #30: "Prototipo de motor de reglas para el horario"
# (added in this order so the shared-string table grows in the same
# sequence as the target workbook: caca, viernes 14, Prototipo de motor..., al generar...)
$ws.Range("I31").Value = "viernes 14"
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = "Prototipo de motor de reglas para el horario"
$ws.Range("D31").Value = "al generar un horario  cumple con los demas requisitos"
$ws.Range("F31").Value = "A"
$ws.Range("G31").Value = 100
$ws.Rows.Item(31).RowHeight = 21

# Update the window scroll / selection to match the edited view
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H36").Select()
